# The "часовой посыл" (hourly message) block occupying rows 8-21 is being
# moved from the 11:xx hour to the 21:xx hour: every row's time-range in
# column B is shifted forward by 10 hours, and the trigger-minute list in
# column C has its "13" entry changed to "12" (row text/content in columns
# A and D is unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTriggers = '["01", "04", "08", "12", "17", "22", "26", "30"]'

$times = @{
    8  = "20:55 - 20:59"
    9  = "21:00 - 21:04"
    10 = "21:05 - 21:09"
    11 = "21:10 - 21:14"
    12 = "21:15 - 21:19"
    13 = "21:20 - 21:24"
    14 = "21:25 - 21:29"
    15 = "21:30 - 21:34"
    16 = "21:35 - 21:39"
    17 = "21:40 - 21:44"
    18 = "21:45 - 21:49"
    19 = "21:50 - 21:54"
    20 = "21:55 - 21:59"
    21 = "22:00 - 22:04"
}

foreach ($row in 8..21) {
    $ws.Cells.Item($row, 2).Value = $times[$row]
    $ws.Cells.Item($row, 3).Value = $newTriggers
}

# Reflect the author's final cursor position (row 8 is now the top of the
# "часовой посыл" block after the edit).
$null = $ws.Range("A8").Select()
